# Update NATMI LR-pair (Fzd8-Ckap4) TPM-derived metrics in rows 2-10
# (columns E-T) to reflect re-run of the pipeline with updated TPM data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 3.235341333333333
$ws.Cells.Item(2, 8).Value = 9.706023999999999
$ws.Cells.Item(2, 9).Value = 0.2153734454473681
$ws.Cells.Item(2, 10).Value = 0.2153734454473681
$ws.Cells.Item(2, 13).Value = 2.198890666666667
$ws.Cells.Item(2, 14).Value = 6.596672
$ws.Cells.Item(2, 15).Value = 0.104354912045016
$ws.Cells.Item(2, 16).Value = 0.104354912045016
$ws.Cells.Item(2, 17).Value = 7.114161861347555
$ws.Cells.Item(2, 18).Value = 64.02745675212799
$ws.Cells.Item(2, 19).Value = 0.02247527695649215
$ws.Cells.Item(2, 20).Value = 0.02247527695649215
# Row 3 (ECs -> FAPs)
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 3.235341333333333
$ws.Cells.Item(3, 8).Value = 9.706023999999999
$ws.Cells.Item(3, 9).Value = 0.2153734454473681
$ws.Cells.Item(3, 10).Value = 0.2153734454473681
$ws.Cells.Item(3, 15).Value = 0.7101322821622501
$ws.Cells.Item(3, 16).Value = 0.7101322821622501
$ws.Cells.Item(3, 17).Value = 48.41167415378666
$ws.Cells.Item(3, 18).Value = 435.7050673840799
$ws.Cells.Item(3, 19).Value = 0.1529436363326864
$ws.Cells.Item(3, 20).Value = 0.1529436363326864
# Row 4 (ECs -> MuSCs)
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 3.235341333333333
$ws.Cells.Item(4, 8).Value = 9.706023999999999
$ws.Cells.Item(4, 9).Value = 0.2153734454473681
$ws.Cells.Item(4, 10).Value = 0.2153734454473681
$ws.Cells.Item(4, 15).Value = 0.1855128057927339
$ws.Cells.Item(4, 16).Value = 0.1855128057927339
$ws.Cells.Item(4, 17).Value = 12.64691907548089
$ws.Cells.Item(4, 18).Value = 113.822271679328
$ws.Cells.Item(4, 19).Value = 0.03995453215818955
$ws.Cells.Item(4, 20).Value = 0.03995453215818956
# Row 5 (FAPs -> ECs)
$ws.Cells.Item(5, 9).Value = 0.4841904166376352
$ws.Cells.Item(5, 10).Value = 0.4841904166376352
$ws.Cells.Item(5, 13).Value = 2.198890666666667
$ws.Cells.Item(5, 14).Value = 6.596672
$ws.Cells.Item(5, 15).Value = 0.104354912045016
$ws.Cells.Item(5, 16).Value = 0.104354912045016
$ws.Cells.Item(5, 17).Value = 15.99365691772445
$ws.Cells.Item(5, 18).Value = 143.94291225952
$ws.Cells.Item(5, 19).Value = 0.05052764834126008
$ws.Cells.Item(5, 20).Value = 0.05052764834126007
# Row 6 (FAPs -> FAPs)
$ws.Cells.Item(6, 9).Value = 0.4841904166376352
$ws.Cells.Item(6, 10).Value = 0.4841904166376352
$ws.Cells.Item(6, 15).Value = 0.7101322821622501
$ws.Cells.Item(6, 16).Value = 0.7101322821622501
$ws.Cells.Item(6, 19).Value = 0.3438392455679746
$ws.Cells.Item(6, 20).Value = 0.3438392455679746
# Row 7 (FAPs -> MuSCs)
$ws.Cells.Item(7, 9).Value = 0.4841904166376352
$ws.Cells.Item(7, 10).Value = 0.4841904166376352
$ws.Cells.Item(7, 15).Value = 0.1855128057927339
$ws.Cells.Item(7, 16).Value = 0.1855128057927339
$ws.Cells.Item(7, 19).Value = 0.08982352272840051
$ws.Cells.Item(7, 20).Value = 0.08982352272840051
# Row 8 (MuSCs -> ECs)
$ws.Cells.Item(8, 9).Value = 0.3004361379149967
$ws.Cells.Item(8, 10).Value = 0.3004361379149967
$ws.Cells.Item(8, 13).Value = 2.198890666666667
$ws.Cells.Item(8, 14).Value = 6.596672
$ws.Cells.Item(8, 15).Value = 0.104354912045016
$ws.Cells.Item(8, 16).Value = 0.104354912045016
$ws.Cells.Item(8, 17).Value = 9.923931474865778
$ws.Cells.Item(8, 18).Value = 89.31538327379199
$ws.Cells.Item(8, 19).Value = 0.03135198674726378
$ws.Cells.Item(8, 20).Value = 0.03135198674726378
# Row 9 (MuSCs -> FAPs)
$ws.Cells.Item(9, 9).Value = 0.3004361379149967
$ws.Cells.Item(9, 10).Value = 0.3004361379149967
$ws.Cells.Item(9, 15).Value = 0.7101322821622501
$ws.Cells.Item(9, 16).Value = 0.7101322821622501
$ws.Cells.Item(9, 18).Value = 607.7887059983699
$ws.Cells.Item(9, 19).Value = 0.2133494002615892
$ws.Cells.Item(9, 20).Value = 0.2133494002615892
# Row 10 (MuSCs -> MuSCs)
$ws.Cells.Item(10, 9).Value = 0.3004361379149967
$ws.Cells.Item(10, 10).Value = 0.3004361379149967
$ws.Cells.Item(10, 15).Value = 0.1855128057927339
$ws.Cells.Item(10, 16).Value = 0.1855128057927339
$ws.Cells.Item(10, 19).Value = 0.05573475090614379
$ws.Cells.Item(10, 20).Value = 0.05573475090614379
